# Update the "dSF" (F) column values for the specific rows that were
# repulled/recalculated, per the commit message "repull data, push all
# data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -5
$ws.Range("F4").Value = -3
$ws.Range("F5").Value = 6
$ws.Range("F6").Value = -2
$ws.Range("F9").Value = -3
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = -5
$ws.Range("F15").Value = -1
$ws.Range("F16").Value = 4
